$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: String Address adr
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# Row 10: Datatype Environment
$ws.Range("B10").Value = "Environment"

# Row 11: import com.example.beans
$ws.Range("B11").Value = "import"
$ws.Range("C11").Value = "com.example.beans"

# Apply style to B10, C10, B11, C11 matching style index 2 (same as B2/B3 etc.)
$ws.Range("B10:C11").Style = $ws.Range("B2").Style

# Column widths (best-fit for new longest content per column)
$ws.Columns.Item(2).ColumnWidth = 15.7109375
$ws.Columns.Item(3).ColumnWidth = 19.140625
$ws.Columns.Item(4).ColumnWidth = 25

# Select C10 as active cell
$ws.Range("C10").Select()
